# Sprint 3 effort update for Travis Thayer (GPTMS)

$wb = $excel.ActiveWorkbook

$wsEvidence = $wb.Worksheets.Item("Evidence")
$wsCount    = $wb.Worksheets.Item("Count")

# ---------------------------------------------------------------------
# "Evidence" sheet — fill in the Sprint # 3 row (row 6) with the new
# effort entry describing the PHP/Apache migration work.
# ---------------------------------------------------------------------
$wsEvidence.Range("C6").Value = 16
$wsEvidence.Range("D6").Value = 2
$wsEvidence.Range("E6").Value = "Switched from Python/Flask to PHP/Apache for easier MVC integration. Created an endpoint for the login page and registration page. Created a User class, user controller, SQL connection via mysqli, and related SQL queries for login and registration funtionality. "
$wsEvidence.Range("F6").Value = "profileController.php; ProfileQueries.php; databaseConnection.php; User.php"
$wsEvidence.Range("G6").Value = "https://github.com/quentinxs/GPTMS/blob/travis/api/User/User.php; https://github.com/quentinxs/GPTMS/blob/travis/api/User/profileController.php; https://github.com/quentinxs/GPTMS/blob/travis/api/User/profileQueries.php"
$wsEvidence.Range("H6").Value = "profileController.php; ProfileQueries.php; databaseConnection.php; User.php"

# match the wrap formatting used by the row above it, and let the row
# grow tall enough to show the new description text
$wsEvidence.Range("E6:H6").WrapText = $true
$wsEvidence.Rows.Item(6).RowHeight = 93.6

# ---------------------------------------------------------------------
# "Count" sheet — record Travis Thayer's objective tally for Sprint 3
# (row 8: Sprint# = 3).
# ---------------------------------------------------------------------
$wsCount.Range("B8").Value = "Travis Thayer"
$wsCount.Range("C8").Value = 0
$wsCount.Range("D8").Value = 1
$wsCount.Range("E8").Value = 0
$wsCount.Range("F8").Value = 0
$wsCount.Range("G8").Value = 0
$wsCount.Range("H8").Value = 0
$wsCount.Range("I8").Value = 0

# ---------------------------------------------------------------------
# Final cursor / selection state: Count was left on D9, and Evidence
# (the active tab) was left on C9.
# ---------------------------------------------------------------------
$wsCount.Range("D9").Select()
$wsEvidence.Activate()
$wsEvidence.Range("C9").Select()
